$wb = $excel.ActiveWorkbook

# --- Content changes: replace the "<PlayaNombre1>" placeholder with the
# concrete literal value "CPA_Playa1" that testers should enter. ---

$wsPre = $wb.Worksheets.Item("Precondiciones")
$wsPre.Range("B3").Value = '"CPA_Playa1" es el nombre de la playa <Playa1>'

$wsPasos = $wb.Worksheets.Item("Pasos")
$wsPasos.Range("B3").Value = 'Ingreso "CPA_Playa1" en el campo nombre de playa'
$wsPasos.Range("C5").Value = 'Se muestra el siguiente mesaje confirmacion "Esta seguro que desea eliminar la playa CPA_Playa1"'
$wsPasos.Range("C6").Value = 'Se muestra un mensaje que diga: "La playa CPA_Playa1 ha sido eliminada con éxito"'

# The shorter replacement text in C5 now wraps to fewer lines, so the row
# shrinks from its previous 3-line height down to the standard 2-line height
# used by the other wrapped rows in this sheet.
$wsPasos.Rows.Item(5).RowHeight = 26.25

# --- View / selection changes: Precondiciones selection moves to B3,
# Pasos selection moves to B11 and becomes the active (displayed) tab,
# Control de cambios loses the active-tab flag but keeps its D3 selection. ---

$wsCambios = $wb.Worksheets.Item("Control de cambios")

$wsPre.Range("B3").Select() | Out-Null
$wsCambios.Range("D3").Select() | Out-Null
$wsPasos.Range("B11").Select() | Out-Null
$wsPasos.Activate() | Out-Null
